$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 88, shifting existing rows 88:202 down to 89:203
$ws.Rows.Item(88).Insert()

# Populate the newly-inserted row 88 with its data (new weekly price entry)
$ws.Range("A88").Value = 3
$ws.Range("B88").Value = "Femacal de La Calera"
$ws.Range("C88").Value = "Coquimbo"
$ws.Range("D88").Value = 44482
$ws.Range("E88").Value = 5
$ws.Range("F88").Value = 100114013
$ws.Range("G88").Value = "Zanahoria"
$ws.Range("H88").Value = "Sin especificar"
$ws.Range("I88").Value = "Primera"
$ws.Range("J88").Value = 450
$ws.Range("K88").Value = 8000
$ws.Range("L88").Value = 8500
$ws.Range("M88").Value = 8256
$ws.Range("N88").Value = "$/saco 20 kilos"
$ws.Range("O88").Value = "Provincia de Quillota"
$ws.Range("P88").Value = 413
$ws.Range("Q88").Value = 20
$ws.Range("R88").Value = "Hortaliza"
